$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-26
$values = @(
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(6, 6),
    @(5, 6),
    @(8, 9),
    @(8, 8),
    @(6, 6),
    @(6, 7),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(4, 4),
    @(5, 6),
    @(5, 6),
    @(1, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
